$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (D) and 1h volume change (E) refresh.
# D-column values that are plain numeric-looking strings must be forced to
# Text format first, otherwise COM auto-converts them to numbers and loses
# formatting (e.g. trailing zeros "1.40" -> 1.4, "0.0000170" -> 1.7E-05).

$ws.Range("D2").Value = "61.354.77"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.376.20"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.33"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "2.377.30"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.52"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("E15").Value = "  +5.38%  "
$ws.Range("D16").Value = "2.806.37"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "61.278.43"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "2.376.79"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.66"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  -8.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.22"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.87"
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("D28").Value = "2.494.00"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "519.98"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").Value = "0.0₃0902"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.70"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.51"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.13"
$ws.Range("E42").Value = "  +6.68%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.35"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.63"
$ws.Range("E45").Value = "  +6.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.78"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.581"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0908"
$ws.Range("E51").Value = "  +1.83%  "
